# Adds "October 14, 2020" (Excel serial 44118) data to the COVID-19 Mexico
# bitacora workbook: one new row on out_vars / dates_dx / dates_sx /
# dates_deaths, one new column on control_obs, and updates the active
# sheet/selection bookkeeping to match.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) out_vars (sheet 1): append row 137
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("out_vars")

$ws1.Range("A136:J136").Copy() | Out-Null
$ws1.Range("A137:J137").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws1.Range("A137").Value2 = 44118
$ws1.Range("B137").Value2 = 829396
$ws1.Range("C137").Value2 = 987457
$ws1.Range("D137").Value2 = 310814
$ws1.Range("E137").Value2 = 84898
$ws1.Range("F137").Value2 = 23.335535739260859
$ws1.Range("G137").Value2 = 193544
$ws1.Range("H137").Value2 = 16748
$ws1.Range("I137").Value2 = 33375
$ws1.Range("J137").Value2 = 2127667

$ws1.Range("A137").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) dates_dx (sheet 2): fill row 135, append row 136
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("dates_dx")

$ws2.Range("A134:O134").Copy() | Out-Null
$ws2.Range("A135:O135").PasteSpecial($xlPasteFormats) | Out-Null
$ws2.Range("A136:O136").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws2.Range("B135").Value2 = 0
$ws2.Range("C135").Value2 = 1
$ws2.Range("D135").Value2 = 0
$ws2.Range("E135").Value2 = 0
$ws2.Range("F135").Value2 = 1
$ws2.Range("G135").Value2 = 0
$ws2.Range("H135").Value2 = 0
$ws2.Range("I135").Value2 = 0
$ws2.Range("J135").Value2 = 1
$ws2.Range("K135").Value2 = 0
$ws2.Range("L135").Value2 = 0
$ws2.Range("M135").Value2 = 1
$ws2.Range("N135").Value2 = 2
$ws2.Range("O135").Value2 = 5

$ws2.Range("A136").Value2 = 44118
$ws2.Range("B136").Value2 = 0
$ws2.Range("C136").Value2 = 1
$ws2.Range("D136").Value2 = 0
$ws2.Range("E136").Value2 = 0
$ws2.Range("F136").Value2 = 0
$ws2.Range("G136").Value2 = 0
$ws2.Range("H136").Value2 = 0
$ws2.Range("I136").Value2 = 0
$ws2.Range("J136").Value2 = 1
$ws2.Range("K136").Value2 = 0
$ws2.Range("L136").Value2 = 0
$ws2.Range("M136").Value2 = 1
$ws2.Range("N136").Value2 = 2
$ws2.Range("O136").Value2 = 5

$ws2.Range("O136").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) dates_sx (sheet 3): fill row 135, append row 136
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("dates_sx")

$ws3.Range("A130:N130").Copy() | Out-Null
$ws3.Range("A135:N135").PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Range("A136:N136").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws3.Range("B135").Value2 = 0
$ws3.Range("C135").Value2 = 1
$ws3.Range("D135").Value2 = 0
$ws3.Range("E135").Value2 = 0
$ws3.Range("F135").Value2 = 0
$ws3.Range("G135").Value2 = 0
$ws3.Range("H135").Value2 = 1
$ws3.Range("I135").Value2 = 0
$ws3.Range("J135").Value2 = 0
$ws3.Range("K135").Value2 = 1
$ws3.Range("L135").Value2 = 0
$ws3.Range("M135").Value2 = 0
$ws3.Range("N135").Value2 = 0

$ws3.Range("A136").Value2 = 44118
$ws3.Range("B136").Value2 = 0
$ws3.Range("C136").Value2 = 1
$ws3.Range("D136").Value2 = 0
$ws3.Range("E136").Value2 = 0
$ws3.Range("F136").Value2 = 0
$ws3.Range("G136").Value2 = 0
$ws3.Range("H136").Value2 = 1
$ws3.Range("I136").Value2 = 0
$ws3.Range("J136").Value2 = 0
$ws3.Range("K136").Value2 = 1
$ws3.Range("L136").Value2 = 0
$ws3.Range("M136").Value2 = 0
$ws3.Range("N136").Value2 = 0

$ws3.Range("N136").Select() | Out-Null

# ---------------------------------------------------------------------
# 4) dates_deaths (sheet 4): fill row 135, append row 136
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("dates_deaths")

$ws4.Range("A134:M134").Copy() | Out-Null
$ws4.Range("A135:M135").PasteSpecial($xlPasteFormats) | Out-Null
$ws4.Range("A136:M136").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$ws4.Range("B135").Value2 = 0
$ws4.Range("C135").Value2 = 0
$ws4.Range("D135").Value2 = 0
$ws4.Range("E135").Value2 = 1
$ws4.Range("F135").Value2 = 1
$ws4.Range("G135").Value2 = 1
$ws4.Range("H135").Value2 = 0
$ws4.Range("I135").Value2 = 2
$ws4.Range("J135").Value2 = 1
$ws4.Range("K135").Value2 = 2
$ws4.Range("L135").Value2 = 1
$ws4.Range("M135").Value2 = 2

$ws4.Range("A136").Value2 = 44118
$ws4.Range("B136").Value2 = 0
$ws4.Range("C136").Value2 = 0
$ws4.Range("D136").Value2 = 0
$ws4.Range("E136").Value2 = 1
$ws4.Range("F136").Value2 = 1
$ws4.Range("G136").Value2 = 1
$ws4.Range("H136").Value2 = 0
$ws4.Range("I136").Value2 = 2
$ws4.Range("J136").Value2 = 1
$ws4.Range("K136").Value2 = 2
$ws4.Range("L136").Value2 = 1
$ws4.Range("M136").Value2 = 2

$ws4.Range("K139").Select() | Out-Null

# ---------------------------------------------------------------------
# 5) control_obs (sheet 5): append column EG (2020-10-14)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("control_obs")

$egRows = @(1,2,3,4,5,6,7,8,10,11,12,13,14,15,16,18)
foreach ($r in $egRows) {
    $ws5.Range("EF$r").Copy() | Out-Null
    $ws5.Range("EG$r").PasteSpecial($xlPasteFormats) | Out-Null
}
$excel.CutCopyMode = $false

$ws5.Range("EG1").Value2 = 44118
$ws5.Range("EG2").Value2 = 7321
$ws5.Range("EG3").Value2 = 7166
$ws5.Range("EG4").Value2 = 7166
$ws5.Range("EG5").Value2 = 7166
$ws5.Range("EG6").Value2 = 7166
$ws5.Range("EG7").Value2 = 6400
$ws5.Range("EG8").Value2 = 9195
$ws5.Range("EG10").Value2 = 276
$ws5.Range("EG11").Value2 = 276
$ws5.Range("EG12").Value2 = 276
$ws5.Range("EG13").Value2 = 276
$ws5.Range("EG14").Value2 = 276
$ws5.Range("EG15").Value2 = 226
$ws5.Range("EG16").Value2 = 288
$ws5.Range("EG18").Value2 = 1689

$ws5.Range("EF20").Copy() | Out-Null
$ws5.Range("EG20").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$ws5.Range("EG20").Formula = "=SUM(EG2:EG18)"

$ws5.Activate() | Out-Null
$ws5.Range("EK15").Select() | Out-Null

Write-Host "edit complete"
